$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the task in row 9 ("View shopping list" / "Complete functionalty to view
# current shopping list for web application") as completed by Destiny, matching
# the pattern used by the other completed rows (7, 13, 15):
#   Actual Time (E) = 1
#   Completed By (F) = same team member as Assigned Team Member (D)
#   Week totals (H, I) = 0
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = $ws.Range("D9").Value2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

# Update the active selection to reflect where the user left off editing.
$ws.Range("B13").Select()

$wb.Save()
